$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata "type" annotations - dimension -> measure for the curated columns
$ws.Range("D2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("H2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("I2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: dim -> medida for the same curated columns
$ws.Range("D3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"

# Row 4: skos:Concept / URI-Provincia -> xsd:int for the same curated columns
$ws.Range("D4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"

# Row 5: remove the now-unused mapping file references for the curated columns
$ws.Range("D5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
